# Actualización automática 2025-10-16 08:30:08
# Applies the updated sales figures across the three sheets of the workbook.

$wb = $excel.ActiveWorkbook

$wsVentasGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO": per-client / per-category amounts ---
$wsVentasGrupo.Range("M4").Value  = 2198.13
$wsVentasGrupo.Range("L5").Value  = 2526.74
$wsVentasGrupo.Range("M5").Value  = 72.06999999999999
$wsVentasGrupo.Range("D15").Value = 457.92
$wsVentasGrupo.Range("H15").Value = 980.1
$wsVentasGrupo.Range("K35").Value = 152.28
$wsVentasGrupo.Range("L35").Value = 1835.01
$wsVentasGrupo.Range("M35").Value = 4381.68
$wsVentasGrupo.Range("K36").Value = 319.68
$wsVentasGrupo.Range("M36").Value = 1275.22
$wsVentasGrupo.Range("K40").Value = 91.37

# Row 55 "X de 53" progress counters
$wsVentasGrupo.Range("D55").Value = "7 de 53"
$wsVentasGrupo.Range("H55").Value = "2 de 53"
$wsVentasGrupo.Range("K55").Value = "4 de 53"
$wsVentasGrupo.Range("L55").Value = "2 de 53"
$wsVentasGrupo.Range("M55").Value = "9 de 53"

# --- Sheet "VENTA MENSUAL": octubre column (F) per client, plus total row ---
$wsVentaMensual.Range("F4").Value  = 3818.12
$wsVentaMensual.Range("F5").Value  = 2598.81
$wsVentaMensual.Range("F15").Value = 1419.49
$wsVentaMensual.Range("F35").Value = 6368.97
$wsVentaMensual.Range("F36").Value = 2849.01
$wsVentaMensual.Range("F40").Value = 91.37
$wsVentaMensual.Range("F59").Value = 28402.42

# --- Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / CUMPLIMIENTO per grupo ---
$wsCumplimiento.Range("D3").Value  = 4648.32
$wsCumplimiento.Range("E3").Value  = 13020.8270988183
$wsCumplimiento.Range("F3").Value  = 0.2630755165488931

$wsCumplimiento.Range("D6").Value  = 1888.2
$wsCumplimiento.Range("E6").Value  = 1019.38368146026
$wsCumplimiento.Range("F6").Value  = 0.6494052130089338

$wsCumplimiento.Range("D10").Value = 969.41
$wsCumplimiento.Range("E10").Value = 2911.66983534392
$wsCumplimiento.Range("F10").Value = 0.2497784227914745

$wsCumplimiento.Range("D11").Value = 4361.75
$wsCumplimiento.Range("E11").Value = 7469.25
$wsCumplimiento.Range("F11").Value = 0.3686712872960866

$wsCumplimiento.Range("D12").Value = 9873.23
$wsCumplimiento.Range("E12").Value = 42789.89
$wsCumplimiento.Range("F12").Value = 0.187479017574348

$wsCumplimiento.Range("D14").Value = 26381.73
$wsCumplimiento.Range("E14").Value = 72634.77661190613
$wsCumplimiento.Range("F14").Value = 0.2664376971347094
